$wb = $excel.ActiveWorkbook

$wsNorite = $wb.Worksheets.Item("Norite")
$wsAtten  = $wb.Worksheets.Item("Attenuation Coefficients")

# --- Append a second "Mean Free Path" table (rows 17-28) to the
#     "Attenuation Coefficients" sheet, mirroring the existing table
#     (rows 3-13) but driven off the "Norite" composition instead of
#     "Shotcrete Density". Row 29 is the totals row.

$elements = @(
    @{ Row = 17; A = "H ";  B = 1;  C = 3.927;                E = 0.99972000000000005; NoriteRow = 2  },
    @{ Row = 18; A = "C ";  B = 12; C = 2.3380000000000001;    E = 0.98839999999999995; NoriteRow = 3  },
    @{ Row = 19; A = "O";   B = 16; C = 2.7559999999999998;    E = 0.99738000000000004; NoriteRow = 4  },
    @{ Row = 20; A = "Na";  B = 23; C = 2.7;                   E = 1;                   NoriteRow = 5  },
    @{ Row = 21; A = "Mg";  B = 24; C = 3.145;                 E = 0.78879999999999995; NoriteRow = 6  },
    @{ Row = 22; A = "Al";  B = 27; C = 2.9660000000000002;    E = 1;                   NoriteRow = 7  },
    @{ Row = 23; A = "Si";  B = 28; C = 2.99;                  E = 0.92191000000000001; NoriteRow = 8  },
    @{ Row = 24; A = "K";   B = 39; C = 2.6539999999999999;    E = 0.93258099999999999; NoriteRow = 9  },
    @{ Row = 25; A = "Ca";  B = 40; C = 2.7010000000000001;    E = 0.96940999999999999; NoriteRow = 10 },
    @{ Row = 26; A = "Mn "; B = 55; C = 2.5880000000000001;    E = 1;                   NoriteRow = 11 },
    @{ Row = 27; A = "Fe";  B = 56; C = 2.5880000000000001;    E = 0.91754000000000002; NoriteRow = 12 },
    @{ Row = 28; A = "Ti "; B = 48; C = 2.661;                 E = 0.73719999999999997; NoriteRow = 13 }
)

foreach ($el in $elements) {
    $r = $el.Row
    $wsAtten.Range("A$r").Value = $el.A
    $wsAtten.Range("B$r").Value = $el.B
    $wsAtten.Range("C$r").Value = $el.C
    $wsAtten.Range("D$r").Formula = "=C$r*1E-24"
    $wsAtten.Range("E$r").Value = $el.E
    $wsAtten.Range("F$r").Formula = "=Norite!H" + $el.NoriteRow
    $wsAtten.Range("G$r").Formula = "=F$r*E$r"
    $wsAtten.Range("H$r").Formula = "=D$r*F$r"
    $wsAtten.Range("I$r").Formula = "=1/H$r"
}

# Totals row
$wsAtten.Range("A29").Value = "Totals"
$wsAtten.Range("F29").Formula = "=SUM(F17:F28)"
$wsAtten.Range("G29").Formula = "=SUM(G18:G28)"
$wsAtten.Range("H29").Formula = "=SUM(H18:H28)"
$wsAtten.Range("I29").Formula = "=1/H29"

# --- View-state bookkeeping: the user ended up on the Attenuation
#     Coefficients sheet (now the active tab) with a selection near the
#     bottom of the new table, while Norite kept a lingering selection
#     from before the tab switch.
$wsNorite.Range("D36").Select()
$wsAtten.Activate()
$wsAtten.Range("I33").Select()
